# Cck-Cckar (OldD7) NATMI ligand-receptor edge table refreshed against new TPM input.
# Existing ECs/FAPs/MuSCs rows (2-7) get updated detection/expression/specificity
# metrics, and a new "Resolving-Mac" sending cluster contributes two more rows (8-9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number -> full A:T row values (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster, then the 16 NATMI numeric metric columns).
$rowData = [ordered]@{
    2 = @("ECs", "Cck", "Cckar", "FAPs", 3, 1, 0.857847, 2.573541, 0.1498657082557423, 0.1498657082557423, 1, 0.3333333333333333, 0.1713833333333333, 0.51415, 0.4272957328601764, 0.4272957328601762, 0.14702067835, 1.32318610515, 0.06403697763974681, 0.06403697763974678)
    3 = @("ECs", "Cck", "Cckar", "MuSCs", 3, 1, 0.857847, 2.573541, 0.1498657082557423, 0.1498657082557423, 2, 0.6666666666666666, 0.229705, 0.6891149999999999, 0.5727042671398237, 0.5727042671398237, 0.197051745135, 1.773465706215, 0.08582873061599554, 0.08582873061599554)
    4 = @("FAPs", "Cck", "Cckar", "FAPs", 3, 1, 2.378641333333333, 7.135924, 0.4155481899527343, 0.4155481899527343, 1, 0.3333333333333333, 0.1713833333333333, 0.51415, 0.4272957328601764, 0.4272957328601762, 0.4076594805111111, 3.6689353246, 0.1775619683645734, 0.1775619683645733)
    5 = @("FAPs", "Cck", "Cckar", "MuSCs", 3, 1, 2.378641333333333, 7.135924, 0.4155481899527343, 0.4155481899527343, 2, 0.6666666666666666, 0.229705, 0.6891149999999999, 0.5727042671398237, 0.5727042671398237, 0.5463858074733332, 4.91747226726, 0.237986221588161, 0.2379862215881609)
    6 = @("MuSCs", "Cck", "Cckar", "FAPs", 3, 1, 2.389512333333334, 7.168537000000001, 0.4174473515916376, 0.4174473515916375, 1, 0.3333333333333333, 0.1713833333333333, 0.51415, 0.4272957328601764, 0.4272957328601762, 0.4095225887277779, 3.68570329855, 0.1783734720288885, 0.1783734720288884)
    7 = @("MuSCs", "Cck", "Cckar", "MuSCs", 3, 1, 2.389512333333334, 7.168537000000001, 0.4174473515916376, 0.4174473515916375, 2, 0.6666666666666666, 0.229705, 0.6891149999999999, 0.5727042671398237, 0.5727042671398237, 0.5488829305283334, 4.939946374754999, 0.2390738795627491, 0.2390738795627491)
    8 = @("Resolving-Mac", "Cck", "Cckar", "FAPs", 1, 0.3333333333333333, 0.09810400000000001, 0.294312, 0.0171387501998857, 0.0171387501998857, 1, 0.3333333333333333, 0.1713833333333333, 0.51415, 0.4272957328601764, 0.4272957328601762, 0.01681339053333333, 0.1513205148, 0.007323314826967654, 0.007323314826967651)
    9 = @("Resolving-Mac", "Cck", "Cckar", "MuSCs", 1, 0.3333333333333333, 0.09810400000000001, 0.294312, 0.0171387501998857, 0.0171387501998857, 2, 0.6666666666666666, 0.229705, 0.6891149999999999, 0.5727042671398237, 0.5727042671398237, 0.02253497932, 0.20281481388, 0.009815435372918046, 0.009815435372918044)
}

foreach ($rowNum in $rowData.Keys) {
    $values = $rowData[$rowNum]
    for ($colIdx = 0; $colIdx -lt $values.Length; $colIdx++) {
        $colLetter = [char](65 + $colIdx)
        $ws.Range("$colLetter$rowNum").Value = $values[$colIdx]
    }
}
